# Update the "Non-white sample" column of the tuition identity model table
# with the revised regression results (r&r update).

$d = $word.ActiveDocument

$replacements = @(
    @("0.250", "0.059"),
    @("[-0.006, 0.506]", "[-0.087, 0.204]"),
    @("-0.036", "-0.145"),
    @("[-0.234, 0.161]", "[-0.281, -0.008]"),
    @("-0.279", "-0.217"),
    @("[-0.465, -0.094]", "[-0.358, -0.076]"),
    @("-0.066", "-0.214"),
    @("[-0.290, 0.157]", "[-0.370, -0.058]"),
    @("0.040", "0.172"),
    @("[-0.224, 0.304]", "[0.004, 0.341]"),
    @("-0.118", "0.013"),
    @("[-0.345, 0.109]", "[-0.132, 0.158]"),
    @("0.203", "0.233"),
    @("[-0.030, 0.436]", "[0.090, 0.376]"),
    @("-0.320", "-0.124"),
    @("[-0.563, -0.077]", "[-0.265, 0.017]"),
    @("107", "208"),
    @("0.17", "0.13"),
    @("0.99", "1.00")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
